# edit.ps1
# Applies the changes described by the commit "fixed typo and added pdf."
# to the ASIM GitHub Download Guide document:
#  1. Title "Power Water Corporation" -> "Power and Water Corporation"
#     (split across three runs, with the _GoBack bookmark moved here)
#  2. Merge "Browse to " + "GitHub" runs (drop spell-check proofErr markers)
#  3. Merge "GitHub" + " (" runs in the following paragraph (drop proofErr)
#  4. Merge the "You will note..." paragraph's six runs into one run
#  5. Remove the old _GoBack bookmark from the "Download an ASIM release" heading
#  6. Merge the "Tar-GZip" runs (drop proofErr markers) in the downloads paragraph
#  7. Remove proofErr markers around "GitHub" in the footer

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rNs = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

# ---------------------------------------------------------------------------
# 1. Title paragraph: "Power Water Corporation" -> "Power and Water Corporation"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000F1DE1" w:rsidRDefault="000F1DE1" w:rsidP="000F1DE1"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Power </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Water Corporation</w:t></w:r></w:p>
'@
$p1.Range.InsertXML($p1xml)

# ---------------------------------------------------------------------------
# 2. "Browse to " + "GitHub" heading -> single run, no proofErr
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$p11xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007E72D8" w:rsidRPr="00043410" w:rsidRDefault="00885E8D" w:rsidP="00043410"><w:pPr><w:pStyle w:val="Heading1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:lastRenderedPageBreak/><w:t>Browse to GitHub</w:t></w:r></w:p>
'@
$p11.Range.InsertXML($p11xml)

# ---------------------------------------------------------------------------
# 3. "GitHub" + " (" -> single run, no proofErr
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$p12xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00885E8D" w:rsidRDefault="00885E8D"><w:r><w:t>GitHub (</w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidRPr="00DF5F46"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>http://www.github.com</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve">) is a hosted source repository that has been chosen to host </w:t></w:r><w:r w:rsidR="004066D9"><w:t>ASIM</w:t></w:r><w:r><w:t xml:space="preserve"> source code and releases.</w:t></w:r></w:p>
'@
$p12.Range.InsertXML($p12xml)

# ---------------------------------------------------------------------------
# 4. "You will note..." paragraph: merge six runs into one
# ---------------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$p16xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00F77F19" w:rsidRDefault="00F77F19" w:rsidP="00F77F19"><w:r><w:t>You will note that multiple versions of ASIM are available for download. It would be usual to download the latest release unless you have a specific requirement for a prior release. The latest release contains all the latest bug fixes and features.</w:t></w:r></w:p>
'@
$p16.Range.InsertXML($p16xml)

# ---------------------------------------------------------------------------
# 5. "Download an ASIM release" heading: remove the _GoBack bookmark
# ---------------------------------------------------------------------------
$p20 = $d.Paragraphs(20)
$p20xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B70279" w:rsidRPr="00043410" w:rsidRDefault="004066D9" w:rsidP="00043410"><w:pPr><w:pStyle w:val="Heading1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:t>Download</w:t></w:r><w:r w:rsidR="00885E8D"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> an</w:t></w:r><w:r w:rsidR="00B70279" w:rsidRPr="00043410"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:t>ASIM</w:t></w:r><w:r w:rsidR="00885E8D"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> release</w:t></w:r></w:p>
'@
$p20.Range.InsertXML($p20xml)

# ---------------------------------------------------------------------------
# 6. "The other two non-highlighted buttons..." paragraph: merge Tar-GZip runs
# ---------------------------------------------------------------------------
$p24 = $d.Paragraphs(24)
$p24xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000F76DE" w:rsidRDefault="004066D9" w:rsidP="00E34F95"><w:r><w:t>The other two</w:t></w:r><w:r w:rsidR="0000211A"><w:t xml:space="preserve"> non-highlighted</w:t></w:r><w:r><w:t xml:space="preserve"> buttons allow you to download a compressed package of source code in two formats: Zip or Tar-GZip. These both contain the same source code files, but are provided in two different formats for convenience &#8211; normally you would download a Zip file for a Microsoft Windows PC and a Tar-GZip file for Linux platforms.</w:t></w:r></w:p>
'@
$p24.Range.InsertXML($p24xml)

# ---------------------------------------------------------------------------
# 7. Footer: remove proofErr markers around "GitHub"
# ---------------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$fp = $footer.Range.Paragraphs(1)
$fpxml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000F1DE1" w:rsidRDefault="004066D9"><w:pPr><w:pStyle w:val="Footer"/></w:pPr><w:r><w:t>ASIM</w:t></w:r><w:r w:rsidR="000F1DE1"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>GitHub</w:t></w:r><w:r w:rsidR="000F1DE1"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Download Guide</w:t></w:r><w:r w:rsidR="000F1DE1"><w:ptab w:relativeTo="margin" w:alignment="center" w:leader="none"/></w:r><w:r w:rsidR="000F1DE1"><w:ptab w:relativeTo="margin" w:alignment="right" w:leader="none"/></w:r><w:r w:rsidR="000F1DE1"><w:t xml:space="preserve">Page </w:t></w:r><w:r w:rsidR="000F1DE1"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="000F1DE1"><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r><w:r w:rsidR="000F1DE1"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="0008450E"><w:rPr><w:noProof/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="000F1DE1"><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@
$fp.Range.InsertXML($fpxml)

Write-Host "Edits applied successfully."
